$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.053018007093146
$ws.Range("D2").Value = 1.061105909166578
$ws.Range("E2").Value = 1.060020846678676
$ws.Range("F2").Value = 1.071333808338326
$ws.Range("I2").Value = 1.052944360494309
$ws.Range("J2").Value = 1.05803755917594
$ws.Range("K2").Value = 1.063830872483262
$ws.Range("L2").Value = 1.062748762451211
$ws.Range("M2").Value = 1.074031259345916
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.053909042993574
$ws.Range("D3").Value = 1.061824267095554
$ws.Range("E3").Value = 1.060798764278584
$ws.Range("F3").Value = 1.072157001239556
$ws.Range("I3").Value = 1.053212037405383
$ws.Range("J3").Value = 1.058579594282444
$ws.Range("K3").Value = 1.064364029088813
$ws.Range("L3").Value = 1.063341115819466
$ws.Range("M3").Value = 1.074670969082837
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.054486227100928
$ws.Range("D4").Value = 1.062289602363471
$ws.Range("E4").Value = 1.061303013420158
$ws.Range("F4").Value = 1.07269053622579
$ws.Range("I4").Value = 1.053384280490474
$ws.Range("J4").Value = 1.058930278682203
$ws.Range("K4").Value = 1.064708852198147
$ws.Range("L4").Value = 1.063724629630746
$ws.Range("M4").Value = 1.075085123189501
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.054729023267185
$ws.Range("D5").Value = 1.062485349629588
$ws.Range("E5").Value = 1.061515209761781
$ws.Range("F5").Value = 1.072915041796815
$ws.Range("I5").Value = 1.053456460439743
$ws.Range("J5").Value = 1.05907769354267
$ws.Range("K5").Value = 1.064853775072621
$ws.Range("L5").Value = 1.063885910519064
$ws.Range("M5").Value = 1.075259284557543
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.054769798445487
$ws.Range("D6").Value = 1.062518223463274
$ws.Range("E6").Value = 1.06155085072549
$ws.Range("F6").Value = 1.072952749383255
$ws.Range("I6").Value = 1.053468566199314
$ws.Range("J6").Value = 1.059102444356297
$ws.Range("K6").Value = 1.064878105833176
$ws.Range("L6").Value = 1.063912993280928
$ws.Range("M6").Value = 1.075288529964186
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.054489470775632
$ws.Range("D7").Value = 1.062292217477373
$ws.Range("E7").Value = 1.061305847975071
$ws.Range("F7").Value = 1.072693535267298
$ws.Range("I7").Value = 1.05338524587132
$ws.Range("J7").Value = 1.058932248498867
$ws.Range("K7").Value = 1.064710788826469
$ws.Range("L7").Value = 1.063726784472503
$ws.Range("M7").Value = 1.075087450143756
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.053319006685156
$ws.Range("D8").Value = 1.061348575047739
$ws.Range("E8").Value = 1.060283563403291
$ws.Range("F8").Value = 1.071611828152147
$ws.Range("I8").Value = 1.053035021824434
$ws.Range("J8").Value = 1.058220751789017
$ws.Range("K8").Value = 1.064011088614155
$ws.Range("L8").Value = 1.062948904347877
$ws.Range("M8").Value = 1.074247405984984
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.051261346774266
$ws.Range("D9").Value = 1.059689734649113
$ws.Range("E9").Value = 1.058489017726863
$ws.Range("F9").Value = 1.069712496700278
$ws.Range("I9").Value = 1.052410551353993
$ws.Range("J9").Value = 1.056966686436869
$ws.Range("K9").Value = 1.062776921874709
$ws.Range("L9").Value = 1.061579937416452
$ws.Range("M9").Value = 1.072768884507635
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.049892921084038
$ws.Range("D10").Value = 1.058586616465713
$ws.Range("E10").Value = 1.057297362969619
$ws.Range("F10").Value = 1.068450938233449
$ws.Range("I10").Value = 1.051989355970345
$ws.Range("J10").Value = 1.056130503079364
$ws.Range("K10").Value = 1.061953410325995
$ws.Range("L10").Value = 1.060668555365912
$ws.Range("M10").Value = 1.07178446386943
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.049301190159279
$ws.Range("D11").Value = 1.05810963362659
$ws.Range("E11").Value = 1.056782501980701
$ws.Range("F11").Value = 1.067905798036585
$ws.Range("I11").Value = 1.051805825725086
$ws.Range("J11").Value = 1.055768408735586
$ws.Range("K11").Value = 1.061596662851529
$ws.Range("L11").Value = 1.06027423340318
$ws.Range("M11").Value = 1.071358516317762
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.049081517574641
$ws.Range("D12").Value = 1.057932564024835
$ws.Range("E12").Value = 1.05659143172121
$ws.Range("F12").Value = 1.067703479233458
$ws.Range("I12").Value = 1.051737482441218
$ws.Range("J12").Value = 1.055633908785934
$ws.Range("K12").Value = 1.061464127974671
$ws.Range("L12").Value = 1.06012781298416
$ws.Range("M12").Value = 1.07120034895368
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.049128632516925
$ws.Range("D13").Value = 1.057970541364264
$ws.Range("E13").Value = 1.05663240912727
$ws.Range("F13").Value = 1.067746869556939
$ws.Range("I13").Value = 1.051752150077551
$ws.Range("J13").Value = 1.055662759551883
$ws.Range("K13").Value = 1.061492558174879
$ws.Range("L13").Value = 1.060159218452737
$ws.Range("M13").Value = 1.071234274170619
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.049283029444817
$ws.Range("D14").Value = 1.058094994885749
$ws.Range("E14").Value = 1.056766704533803
$ws.Range("F14").Value = 1.0678890708117
$ws.Range("I14").Value = 1.051800179951751
$ws.Range("J14").Value = 1.05575729096122
$ws.Range("K14").Value = 1.061585707938605
$ws.Range("L14").Value = 1.060262129252137
$ws.Range("M14").Value = 1.071345441155684
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.049378174781255
$ws.Range("D15").Value = 1.058171688515695
$ws.Range("E15").Value = 1.056849471215075
$ws.Range("F15").Value = 1.067976708350671
$ws.Range("I15").Value = 1.051829749976759
$ws.Range("J15").Value = 1.055815534655861
$ws.Range("K15").Value = 1.061643097575772
$ws.Range("L15").Value = 1.0603255424353
$ws.Range("M15").Value = 1.071413941278208
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.049932209432676
$ws.Range("D16").Value = 1.058618286605096
$ws.Range("E16").Value = 1.057331556626968
$ws.Range("F16").Value = 1.068487141208185
$ws.Range("I16").Value = 1.052001512108056
$ws.Range("J16").Value = 1.056154533742001
$ws.Range("K16").Value = 1.061977083187496
$ws.Range("L16").Value = 1.060694731916308
$ws.Range("M16").Value = 1.071812739328616
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.050279957673055
$ws.Range("D17").Value = 1.058898607607761
$ws.Range("E17").Value = 1.057634260749678
$ws.Range("F17").Value = 1.068807624189348
$ws.Range("I17").Value = 1.052108946650097
$ws.Range("J17").Value = 1.056367174043776
$ws.Range("K17").Value = 1.062186541198804
$ws.Range("L17").Value = 1.060926399349553
$ws.Range("M17").Value = 1.072062979678403
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.050482871048108
$ws.Range("D18").Value = 1.059062179081822
$ws.Range("E18").Value = 1.057810932180458
$ws.Range("F18").Value = 1.068994664715341
$ws.Range("I18").Value = 1.052171500403442
$ws.Range("J18").Value = 1.056491201383666
$ws.Range("K18").Value = 1.06230869884035
$ws.Range("L18").Value = 1.061061557101785
$ws.Range("M18").Value = 1.072208970748848
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.050552072401511
$ws.Range("D19").Value = 1.059117963696432
$ws.Range("E19").Value = 1.057871191047818
$ws.Range("F19").Value = 1.069058459002381
$ws.Range("I19").Value = 1.052192810766175
$ws.Range("J19").Value = 1.056533491099186
$ws.Range("K19").Value = 1.062350348741722
$ws.Range("L19").Value = 1.06110764745922
$ws.Range("M19").Value = 1.072258754992814
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.050242639529398
$ws.Range("D20").Value = 1.058868525084009
$ws.Range("E20").Value = 1.057601772150455
$ws.Range("F20").Value = 1.068773228191047
$ws.Range("I20").Value = 1.052097431415872
$ws.Range("J20").Value = 1.056344359965349
$ws.Range("K20").Value = 1.062164069948277
$ws.Range("L20").Value = 1.060901540507531
$ws.Range("M20").Value = 1.072036128132173
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.049237560000839
$ws.Range("D21").Value = 1.058058343566846
$ws.Range("E21").Value = 1.056727153120657
$ws.Range("F21").Value = 1.067847191364069
$ws.Range("I21").Value = 1.051786041092141
$ws.Range("J21").Value = 1.055729453862288
$ws.Range("K21").Value = 1.061558278278909
$ws.Range("L21").Value = 1.060231823241971
$ws.Range("M21").Value = 1.071312703897976
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.048606336826702
$ws.Range("D22").Value = 1.057549547423537
$ws.Range("E22").Value = 1.056178241639776
$ws.Range("F22").Value = 1.067265942843984
$ws.Range("I22").Value = 1.051589263004012
$ws.Range("J22").Value = 1.055342827097657
$ws.Range("K22").Value = 1.061177260692187
$ws.Range("L22").Value = 1.059811025630284
$ws.Range("M22").Value = 1.070858139837645
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.048940892361828
$ws.Range("D23").Value = 1.05781921271404
$ws.Range("E23").Value = 1.056469134984968
$ws.Range("F23").Value = 1.067573979459658
$ws.Range("I23").Value = 1.051693672811919
$ws.Range("J23").Value = 1.055547785881644
$ws.Range("K23").Value = 1.061379257435859
$ws.Range("L23").Value = 1.060034071387828
$ws.Range("M23").Value = 1.071099085670718
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.050259501752225
$ws.Range("D24").Value = 1.05888211788105
$ws.Range("E24").Value = 1.057616452012681
$ws.Range("F24").Value = 1.068788769927947
$ws.Range("I24").Value = 1.052102634997599
$ws.Range("J24").Value = 1.056354668672157
$ws.Range("K24").Value = 1.062174223787703
$ws.Range("L24").Value = 1.060912773055494
$ws.Range("M24").Value = 1.072048261097251
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.051792717642588
$ws.Range("D25").Value = 1.06011810284193
$ws.Range("E25").Value = 1.058952128357274
$ws.Range("F25").Value = 1.070202705485457
$ws.Range("I25").Value = 1.052572855841491
$ws.Range("J25").Value = 1.057290922414745
$ws.Range("K25").Value = 1.063096118511559
$ws.Range("L25").Value = 1.061933631584943
$ws.Range("M25").Value = 1.073150902049936
